$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("N2").Value = 0.1094510731625988
$ws.Range("N3").Value = 0.001000000000000013
$ws.Range("N4").Value = 0.001000000000000013
$ws.Range("N5").Value = 0.02061320386266954
$ws.Range("N6").Value = 0.001
$ws.Range("N7").Value = 0.001000000000000009
$ws.Range("N9").Value = 0.001000000000000006
$ws.Range("N10").Value = 0.15
$ws.Range("N11").Value = 0.05499632803248732
$ws.Range("N12").Value = 0.15
$ws.Range("N14").Value = 0.001
$ws.Range("N15").Value = 0.15
$ws.Range("N16").Value = 0.1207720526009886
$ws.Range("N17").Value = 0.001000000000000001
$ws.Range("N18").Value = 0.001
$ws.Range("N19").Value = 0.15
$ws.Range("N20").Value = 0.07382578160983198
$ws.Range("N21").Value = 0.01034156073142385
$ws.Range("B22").Value = -0.004627277821884304
$ws.Range("C22").Value = 0.02530191884529346
$ws.Range("D22").Value = 0.01936141600946563
$ws.Range("E22").Value = 0.01704807988479111
$ws.Range("F22").Value = 0.02968661994640601
$ws.Range("G22").Value = 0.03730654167890406
$ws.Range("H22").Value = -0.01424328450558998
$ws.Range("I22").Value = 0.01711436841469424
$ws.Range("J22").Value = 0.02158004386232901
$ws.Range("K22").Value = 0.04355630760117607
$ws.Range("L22").Value = 0.0008706594484029512
$ws.Range("M22").Value = 0.01176136722928472
$ws.Range("B23").Value = 0.9953834115342403
$ws.Range("C23").Value = 1.025624729217677
$ws.Range("D23").Value = 1.019550063753326
$ws.Range("E23").Value = 1.017194227730828
$ws.Range("F23").Value = 1.030131660650109
$ws.Range("G23").Value = 1.0380111657588
$ws.Range("H23").Value = 0.9858576711892266
$ws.Range("I23").Value = 1.017261658275721
$ws.Range("J23").Value = 1.021814577049407
$ws.Range("K23").Value = 1.044518807004616
$ws.Range("L23").Value = 1.000871038582365
$ws.Range("M23").Value = 1.011830804065861
$ws.Range("N23").Value = 1.133808273796276
